$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# B8 ("mailFilterTypeSubject" value) gains a trailing space: "SUBJECT" -> "SUBJECT "
$ws.Range("B8").Value = "SUBJECT "

# B2 (emaderagheb@gmail.com) is re-entered with its quote-prefix (apostrophe) so the
# cell keeps its quotePrefix formatting but is normalized onto the plain quotePrefix
# style (drops the stray applyFill variant of that style).
$ws.Range("B2").Value = "'emaderagheb@gmail.com"

# Cursor/selection ends on B2.
$ws.Range("B2").Select() | Out-Null
